$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 3649.3333
$ws.Range("I18").Value = 3649.3333
$ws.Range("K18").Value = 3649.3333
$ws.Range("M18").Value = -3365.3333

$ws.Range("H39").Value = 588.0909
$ws.Range("I39").Value = 430
$ws.Range("J39").Value = 1299.5
$ws.Range("K39").Value = 1290
$ws.Range("L39").Value = 3898.5
$ws.Range("M39").Value = -994
$ws.Range("N39").Value = -4490.5

$ws.Range("H100").Value = 9050.833000000001
$ws.Range("I100").Value = 6950
$ws.Range("J100").Value = 10101.25
$ws.Range("K100").Value = 6950
$ws.Range("L100").Value = 10101.25
$ws.Range("M100").Value = -6409
$ws.Range("N100").Value = -11183.25

$ws.Range("H129").Value = 2621.3333
$ws.Range("I129").Value = 2297.4
$ws.Range("J129").Value = 2852.7144
$ws.Range("K129").Value = 6892.200000000001
$ws.Range("L129").Value = 8558.143199999999
$ws.Range("M129").Value = -1892.200000000001
$ws.Range("N129").Value = -18558.1432

$ws.Range("H132").Value = 1609.7441
$ws.Range("I132").Value = 1670.0555
$ws.Range("J132").Value = 1299.5714
$ws.Range("K132").Value = 5010.166499999999
$ws.Range("L132").Value = 3898.7142
$ws.Range("M132").Value = -2480.166499999999
$ws.Range("N132").Value = -8958.7142

$ws.Range("H135").Value = 1316.0358
$ws.Range("I135").Value = 1080.5491
$ws.Range("J135").Value = 3718
$ws.Range("K135").Value = 9724.9419
$ws.Range("L135").Value = 33462
$ws.Range("M135").Value = -7189.9419
$ws.Range("N135").Value = -38532

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2353.942
$ws.Range("I32").Value = 2231.4285
$ws.Range("J32").Value = 7499.5
$ws.Range("K32").Value = 2231.4285
$ws.Range("L32").Value = 7499.5
$ws.Range("M32").Value = -1944.4285
$ws.Range("N32").Value = -8073.5

$ws.Range("H61").Value = 3978.2258
$ws.Range("I61").Value = 1946.8214
$ws.Range("K61").Value = 1946.8214
$ws.Range("M61").Value = -1734.8214

$ws.Range("H74").Value = 10755792
$ws.Range("I74").Value = 11907098
$ws.Range("J74").Value = 10263
$ws.Range("K74").Value = 11907098
$ws.Range("L74").Value = 10263
$ws.Range("M74").Value = -11906224
$ws.Range("N74").Value = -12011

$ws.Range("H77").Value = 10755792
$ws.Range("I77").Value = 11907098
$ws.Range("J77").Value = 10263
$ws.Range("K77").Value = 59535490
$ws.Range("L77").Value = 51315
$ws.Range("M77").Value = -59531122
$ws.Range("N77").Value = -60051

$ws.Range("H110").Value = 6503.1113
$ws.Range("I110").Value = 3013.5715
$ws.Range("J110").Value = 18716.5
$ws.Range("K110").Value = 3013.5715
$ws.Range("L110").Value = 18716.5
$ws.Range("M110").Value = -968.5715
$ws.Range("N110").Value = -22806.5

$ws.Range("H132").Value = 4561.4053
$ws.Range("I132").Value = 3797.6365
$ws.Range("J132").Value = 10862.5
$ws.Range("K132").Value = 11392.9095
$ws.Range("L132").Value = 32587.5
$ws.Range("M132").Value = -8862.9095
$ws.Range("N132").Value = -37647.5

$ws.Range("H136").Value = 3978.2258
$ws.Range("I136").Value = 1946.8214
$ws.Range("K136").Value = 5840.4642
$ws.Range("M136").Value = -3290.4642

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2127.75
$ws.Range("I107").Value = 1005.5
$ws.Range("J107").Value = 3250
$ws.Range("K107").Value = 1005.5
$ws.Range("L107").Value = 3250
$ws.Range("M107").Value = 914.5
$ws.Range("N107").Value = -7090

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2586.842
$ws.Range("I58").Value = 1142.5625
$ws.Range("J58").Value = 10289.667
$ws.Range("K58").Value = 1142.5625
$ws.Range("L58").Value = 10289.667
$ws.Range("M58").Value = -939.5625
$ws.Range("N58").Value = -10695.667

$ws.Range("H97").Value = 99998.5
$ws.Range("J97").Value = 99998.5
$ws.Range("L97").Value = 99998.5
$ws.Range("N97").Value = -101980.5

$ws.Range("H99").Value = 2384.7693
$ws.Range("I99").Value = 1688.5
$ws.Range("J99").Value = 3498.8
$ws.Range("K99").Value = 1688.5
$ws.Range("L99").Value = 3498.8
$ws.Range("M99").Value = -190.5
$ws.Range("N99").Value = -6494.8

$ws.Range("H105").Value = 3146.25
$ws.Range("I105").Value = 879.8570999999999
$ws.Range("J105").Value = 19011
$ws.Range("K105").Value = 879.8570999999999
$ws.Range("L105").Value = 19011
$ws.Range("M105").Value = 867.1429000000001
$ws.Range("N105").Value = -22505

$ws.Range("H107").Value = 1364.8096
$ws.Range("I107").Value = 1231.8
$ws.Range("J107").Value = 1697.3334
$ws.Range("K107").Value = 1231.8
$ws.Range("L107").Value = 1697.3334
$ws.Range("M107").Value = 688.2
$ws.Range("N107").Value = -5537.3334

$ws.Range("H126").Value = 2384.7693
$ws.Range("I126").Value = 1688.5
$ws.Range("J126").Value = 3498.8
$ws.Range("K126").Value = 5065.5
$ws.Range("L126").Value = 10496.4
$ws.Range("M126").Value = -2595.5
$ws.Range("N126").Value = -15436.4

$ws.Range("H132").Value = 3309.9333
$ws.Range("I132").Value = 2869.963
$ws.Range("J132").Value = 7269.6665
$ws.Range("K132").Value = 8609.889000000001
$ws.Range("L132").Value = 21808.9995
$ws.Range("M132").Value = -6079.889000000001
$ws.Range("N132").Value = -26868.9995

$ws.Range("H134").Value = 2479.4883
$ws.Range("I134").Value = 1496.16
$ws.Range("K134").Value = 4488.48
$ws.Range("M134").Value = -1953.48

$ws.Range("H136").Value = 2586.842
$ws.Range("I136").Value = 1142.5625
$ws.Range("J136").Value = 10289.667
$ws.Range("K136").Value = 3427.6875
$ws.Range("L136").Value = 30869.001
$ws.Range("M136").Value = -877.6875
$ws.Range("N136").Value = -35969.001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 2843241.8
$ws.Range("I107").Value = 2019.5
$ws.Range("J107").Value = 10419835
$ws.Range("K107").Value = 6058.5
$ws.Range("L107").Value = 31259505
$ws.Range("M107").Value = -4138.5
$ws.Range("N107").Value = -31263345

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 18786.5
$ws.Range("I46").Value = 8366.666999999999
$ws.Range("J46").Value = 50046
$ws.Range("K46").Value = 8366.666999999999
$ws.Range("L46").Value = 50046
$ws.Range("M46").Value = -8210.666999999999
$ws.Range("N46").Value = -50358

$ws.Range("H80").Value = 5119.3335
$ws.Range("I80").Value = 2206.6667
$ws.Range("J80").Value = 8032
$ws.Range("K80").Value = 2206.6667
$ws.Range("L80").Value = 8032
$ws.Range("M80").Value = -1208.6667
$ws.Range("N80").Value = -10028

$ws.Range("H83").Value = 5119.3335
$ws.Range("I83").Value = 2206.6667
$ws.Range("J83").Value = 8032
$ws.Range("K83").Value = 11033.3335
$ws.Range("L83").Value = 40160
$ws.Range("M83").Value = -6041.333500000001
$ws.Range("N83").Value = -50144

$ws.Range("H113").Value = 2642.3809
$ws.Range("I113").Value = 1927.9286
$ws.Range("J113").Value = 4071.2856
$ws.Range("K113").Value = 1927.9286
$ws.Range("L113").Value = 4071.2856
$ws.Range("M113").Value = 242.0714
$ws.Range("N113").Value = -8411.285599999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4237.4517
$ws.Range("I61").Value = 3073.6667
$ws.Range("J61").Value = 8227.571
$ws.Range("K61").Value = 3073.6667
$ws.Range("L61").Value = 8227.571
$ws.Range("M61").Value = -2871.6667
$ws.Range("N61").Value = -8631.571

$ws.Range("H113").Value = 4237.4517
$ws.Range("I113").Value = 3073.6667
$ws.Range("J113").Value = 8227.571
$ws.Range("K113").Value = 3073.6667
$ws.Range("L113").Value = 8227.571
$ws.Range("M113").Value = -903.6667000000002
$ws.Range("N113").Value = -12567.571

$ws.Range("H136").Value = 4345.282
$ws.Range("I136").Value = 2023.4849
$ws.Range("J136").Value = 17115.166
$ws.Range("K136").Value = 6070.4547
$ws.Range("L136").Value = 51345.49800000001
$ws.Range("M136").Value = -3520.4547
$ws.Range("N136").Value = -56445.49800000001
